$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column G (estoque_atualizado) values per cronjob refresh
$ws.Range("G2").Value = 177
$ws.Range("G3").Value = -40
$ws.Range("G4").Value = -21
$ws.Range("G5").Value = 15
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = -4
$ws.Range("G9").Value = -8
